$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45171 -> 45172) for every
# data row (rows 2 through 406). Update each cell's underlying value while
# preserving its existing date formatting/style.
for ($r = 2; $r -le 406; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
